# "adding averages and more checks"
#
# Training Dashboard sheet, row 3 (the "Endangered by Electricity A safety
# Training (SOPs)" record): the PERIOD TO EXPIRE / LAST UPDATE figures were
# recomputed -
#   H3 (PERIOD TO EXPIRE): -42  -> -50
#   I3 (LAST UPDATE)      : 08-Sep-2025 -> 16-Sep-2025

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# H3 is a plain number - straightforward assignment.
$ws.Range("H3").Value = -50

# I3 holds a literal text label ("16-Sep-2025"), not a real date serial.
# A bare Range.Value assignment of a date-shaped string gets auto-converted
# to a date value (and pulls in a new date number format), so force the
# cell to text first, then restore its original ("General") look by
# re-applying H3's format onto it - this keeps I3 on the same cell style
# it started with while only swapping the displayed text.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "16-Sep-2025"
$ws.Range("H3").Copy()
$ws.Range("I3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
